$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; everything currently at row 20 onward
# shifts down by one (old row 20 -> new row 21, ..., old row 58 -> new row 59).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly entry.
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = 44519
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 100112021
$ws.Cells.Item(20, 7).Value = "Ají"
$ws.Cells.Item(20, 8).Value = "Inferno"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 40
$ws.Cells.Item(20, 11).Value = 17000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 13).Value = 17500
$ws.Cells.Item(20, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(20, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20, 16).Value = 1458
$ws.Cells.Item(20, 17).Value = 12
$ws.Cells.Item(20, 18).Value = "Hortaliza"
